$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '30.654.10'
$ws.Range('E2').Value = '  +1.45%  '
$ws.Range('D3').Value = '1.891.60'
$ws.Range('E3').Value = '  +1.99%  '
Set-TextValue $ws.Range('D4') '1.001'
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue $ws.Range('D5') '239.18'
$ws.Range('E5').Value = '  +1.51%  '
Set-TextValue $ws.Range('D6') '1.002'
$ws.Range('E6').Value = '  +0.09%  '
Set-TextValue $ws.Range('D7') '0.4828'
$ws.Range('E7').Value = '  +0.96%  '
Set-TextValue $ws.Range('D8') '0.2859'
$ws.Range('E8').Value = '  +1.93%  '
Set-TextValue $ws.Range('D9') '0.06552'
$ws.Range('E9').Value = '  +1.06%  '
$ws.Range('D10').Value = '1.881.90'
$ws.Range('E10').Value = '  +1.39%  '
Set-TextValue $ws.Range('D11') '0.07476'
$ws.Range('E11').Value = '  +1.77%  '
$ws.Range('E12').Value = '  +2.95%  '
Set-TextValue $ws.Range('D13') '5.106'
$ws.Range('E13').Value = '  +0.05%  '
Set-TextValue $ws.Range('D14') '88.16'
$ws.Range('E14').Value = '  +1.18%  '
Set-TextValue $ws.Range('D15') '0.6675'
$ws.Range('E15').Value = '  +3.17%  '
$ws.Range('D16').Value = '30.629.14'
$ws.Range('E16').Value = '  +1.49%  '
Set-TextValue $ws.Range('D17') '13.30'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('E18').Value = '  +0.05%  '
Set-TextValue $ws.Range('D19') '0.000007582'
$ws.Range('E19').Value = '  -0.45%  '
Set-TextValue $ws.Range('D20') '232.14'
$ws.Range('E20').Value = '  +3.36%  '
$ws.Range('D21').Value = '2.167.15'
$ws.Range('E21').Value = '  +3.20%  '
$ws.Range('E22').Value = '  -0.02%  '
Set-TextValue $ws.Range('D23') '5.281'
$ws.Range('E23').Value = '  -0.05%  '
Set-TextValue $ws.Range('D24') '6.229'
$ws.Range('E24').Value = '  +2.73%  '
Set-TextValue $ws.Range('D25') '169.46'
$ws.Range('E25').Value = '  +3.28%  '
Set-TextValue $ws.Range('D26') '9.356'
$ws.Range('E26').Value = '  +1.30%  '
Set-TextValue $ws.Range('D27') '18.80'
$ws.Range('E27').Value = '  +2.09%  '
$ws.Range('E28').Value = '  +2.36%  '
$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D29') '0.1024'
$ws.Range('E29').Value = '  +11.35%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D30') '1.401'
$ws.Range('E30').Value = '  -2.95%  '
Set-TextValue $ws.Range('D31') '4.332'
$ws.Range('E31').Value = '  +2.24%  '
Set-TextValue $ws.Range('D32') '4.027'
$ws.Range('E32').Value = '  +1.68%  '
Set-TextValue $ws.Range('D33') '0.05080'
$ws.Range('E33').Value = '  +1.42%  '
Set-TextValue $ws.Range('D34') '1.214'
$ws.Range('E34').Value = '  +5.83%  '
Set-TextValue $ws.Range('D35') '0.7574'
$ws.Range('E35').Value = '  +2.87%  '
Set-TextValue $ws.Range('D36') '2.710'
$ws.Range('E36').Value = '  +0.80%  '
$ws.Range('E37').Value = '  +2.66%  '
Set-TextValue $ws.Range('D38') '2.652'
$ws.Range('E38').Value = '  +1.80%  '
Set-TextValue $ws.Range('D39') '0.9212'
$ws.Range('E39').Value = '  +2.18%  '
Set-TextValue $ws.Range('D40') '2.071'
$ws.Range('E40').Value = '  +0.75%  '
Set-TextValue $ws.Range('D41') '107.18'
$ws.Range('E41').Value = '  +0.85%  '
Set-TextValue $ws.Range('D42') '0.4305'
$ws.Range('E42').Value = '  +1.44%  '
Set-TextValue $ws.Range('D43') '1.003'
$ws.Range('E43').Value = '  +0.25%  '
Set-TextValue $ws.Range('D44') '5.662'
$ws.Range('E44').Value = '  -4.77%  '
Set-TextValue $ws.Range('D45') '7.436'
$ws.Range('E45').Value = '  +0.96%  '
Set-TextValue $ws.Range('D46') '64.19'
$ws.Range('E46').Value = '  -0.21%  '
Set-TextValue $ws.Range('D47') '0.1274'
$ws.Range('E47').Value = '  -3.33%  '
Set-TextValue $ws.Range('D48') '1.491'
$ws.Range('E48').Value = '  -3.00%  '
Set-TextValue $ws.Range('D49') '8.957'
$ws.Range('E49').Value = '  +2.13%  '
Set-TextValue $ws.Range('D50') '33.93'
$ws.Range('E50').Value = '  -0.69%  '
$ws.Range('E51').Value = '  +0.05%  '
